# Apply the updated values to the TDOC balance-sheet worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TDOC")

# Row 4 (Inventory)
$ws.Range("B4").Value = 58000000.0
$ws.Range("C4").Value = 56000000.0

# Row 12 (Accounts Payable)
$ws.Range("B12").Value = 34000000.0
$ws.Range("C12").Value = 46000000.0
$ws.Range("D12").Value = 21000000.0
$ws.Range("E12").Value = 11000000.0
$ws.Range("F12").Value = 8000000.0

# Row 22 (Long Term Tax Liability (Deferred))
$ws.Range("B22").Value = 85000000.0
$ws.Range("C22").Value = 102000000.0
$ws.Range("D22").Value = 18000000.0
$ws.Range("E22").Value = 19000000.0
$ws.Range("F22").Value = 19000000.0

# Row 39 (Net Debt)
$ws.Range("G39").Value = -71654000.0

# Row 40 (Total Debt)
$ws.Range("G40").Value = 445410000.0
